$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheetId=1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 19305
$ws1.Range("F6").Value = 296
$ws1.Range("F9").Value = 7313
$ws1.Range("F13").Value = 28
$ws1.Range("F25").Value = 0
$ws1.Range("F30").Value = 5220
$ws1.Range("F33").Value = 137
$ws1.Range("F36").Value = 12411
$ws1.Range("F45").Value = 94

# Sheet "全部类型" (sheetId=4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F4").Value = 19305
$ws4.Range("F6").Value = 296
$ws4.Range("F9").Value = 7313
$ws4.Range("F13").Value = 28
$ws4.Range("F30").Value = 5220
$ws4.Range("F31").Value = 0
$ws4.Range("F35").Value = 137
$ws4.Range("F38").Value = 12411
$ws4.Range("F41").Value = 3
$ws4.Range("F47").Value = 94
